$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    83.35874938964844,
    86.41961669921875,
    92.64990997314453,
    86.445068359375,
    88.36775207519531,
    88.94989013671875,
    91.51444244384766,
    89.83090209960938,
    90.57490539550781,
    88.1199951171875,
    85.32036590576172,
    88.16629028320312,
    89.97152709960938,
    97.57882690429688,
    101.703010559082,
    130.7038269042969,
    119.657829284668,
    130.9811401367188,
    122.4703216552734,
    125.109733581543,
    128.1146545410156,
    124.8055953979492,
    128.6350860595703,
    127.1221694946289,
    128.81396484375,
    127.8924179077148,
    129.3748626708984,
    119.1679153442383,
    124.4697036743164,
    128.1420135498047,
    134.3809509277344,
    155.2132110595703,
    151.1870422363281,
    185.3318481445312,
    193.1765594482422,
    197.3113403320312,
    180.9517211914062,
    164.3877410888672,
    144.8467864990234,
    129.5536651611328,
    124.5020446777344,
    105.2759170532227,
    109.5275650024414,
    92.63798522949219,
    94.91640472412109,
    80.40143585205078,
    94.48168182373047,
    85.59923553466797
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
